# data-element-mapping.xlsx cleanup
# The source sheet had several cells whose text was corrupted by a stray
# markdown/HTML wrapper (`<span class="bg-success" markdown="1">...</span><!-- new-content -->`)
# left over from a documentation build step, plus a handful of cells where
# the whole wrapper had been typed as manually-colored rich text runs
# (syntax-highlighted HTML, font "Menlo") instead of plain text.
# This script strips the wrapper/markup and replaces each affected cell
# with the clean plain-text value, and normalizes the rich-text cells
# (G9, G10, G17) back to plain, default-styled text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# smart quote characters used verbatim in several of the source strings
$lq = [char]0x201C   # “
$rq = [char]0x201D   # ”

# --- Cells that had the <span ...>...</span><!-- new-content --> wrapper
#     around otherwise-plain text: strip the wrapper, keep the inner text.

$ws.Range("D2").Value = "Loop: 2000E - Patient Event Level or Loop: 2000F Service Level Segment: TRN02 Notes: TRN01 = " + $lq + "1" + $rq + ": Payer Supplied TRN01 = " + $lq + "2" + $rq + " : Provider Supplied (echoed back)"

$ws.Range("C9").Value = "Note that this is an indirect mapping:  Loop 2220D (Service Line Information) Segment SVC01 or SVC04 codes reference the claim service item, which contains the line item."

$ws.Range("D9").Value = "Loop: 2000F Segment: HL01"

$ws.Range("B10").Value = "Attachment Code"

$ws.Range("D10").Value = "Loop: 2000E - Patient Event Level or Loop: 2000F Service Level Segment: HI (LOINC) or Segment: PWK01 Report Type Codes "

$ws.Range("H10").Value = " LOINC Attachment Code.  For prior authorization, [X12] PWK01 Report Type Codes may also be used."
$ws.Range("I10").Value = " LOINC Attachment Code.  For prior authorization, [X12] PWK01 Report Type Codes may also be used."

# --- Cells that were rich text (colored, Menlo-font syntax highlighting of
#     the same wrapper) rather than a plain string: replace with clean
#     plain text and reset formatting back to the default cell style.

$ws.Range("G9").Value = $lq + "AttachmentsNeeded" + $rq + " Task.input.extension"
$ws.Range("G9").Style = "Normal"

$ws.Range("G10").Value = $lq + "AttachmentsNeeded" + $rq + " Task.input"
$ws.Range("G10").Style = "Normal"

$ws.Range("G17").Value = $lq + "QuestionnairesNeeded" + $rq + " Task.input"
$ws.Range("G17").Style = "Normal"
